$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings that look numeric (e.g. "218.27"); Excel
# auto-converts these to Number on plain assignment. Force Text storage by
# temporarily applying a text number format, then clear the format again so
# the cell style matches the untouched cells (style index reverts to General).

$ws.Range('D2').Value = '26.339.53'
$ws.Range('E2').Value = '  +0.21%  '
$ws.Range('D3').Value = '1.686.49'
$ws.Range('E3').Value = '  -0.62%  '
$ws.Range('E4').Value = '  +0.84%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '218.27'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -0.63%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5409'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +2.69%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.012'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  +0.80%  '
$ws.Range('E8').Value = '  +0.87%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06441'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -0.44%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '21.93'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -1.36%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07684'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +2.77%  '
$ws.Range('D12').Value = '1.702.69'
$ws.Range('E12').Value = '  +0.44%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.526'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -0.85%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.5802'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -1.54%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.000008344'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -3.40%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '64.98'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +0.00%  '
$ws.Range('D17').Value = '26.415.22'
$ws.Range('E17').Value = '  +0.21%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '4.928'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -1.46%  '
$ws.Range('E19').Value = '  +0.73%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '10.96'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +0.98%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '190.91'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -0.45%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.231'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -0.67%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.013'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +0.74%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '149.74'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +2.90%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1306'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +4.78%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.846'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +1.83%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '15.65'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -1.89%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.06324'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -7.05%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.406'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +4.47%  '
$ws.Range('E30').Value = '  -0.26%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.581'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -0.98%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.568'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +0.02%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.682'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +0.51%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.038'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +0.34%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.6128'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -1.73%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.416'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +1.46%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.712'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -0.29%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '6.271'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -0.63%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01624'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +0.00%  '
$ws.Range('D40').Value = '1.107.40'
$ws.Range('E40').Value = '  -0.09%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.8815'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +0.27%  '
$ws.Range('E42').Value = '  +0.03%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '101.30'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +0.22%  '
$ws.Range('D44').Value = '1.837.60'
$ws.Range('E44').Value = '  -0.37%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.00000000109'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -1.31%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '57.29'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +0.11%  '
$ws.Range('B47').Value = 'Frax'
$ws.Range('C47').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.013'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +0.42%  '
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '8.195'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +0.11%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.05268'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +0.10%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.4310'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +0.46%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '6.026'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -0.39%  '
